$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("B2").Value = "NSE:AIROLAM"
$ws.Range("C2").Value = "NSE:AARTIPHARM"
$ws.Range("E2").Value = "NSE:IDEA"
$ws.Range("F2").Value = "NSE:AXISBANK"

# --- Row 3 ---
$ws.Range("B3").Value = "NSE:GENUSPAPER"
$ws.Range("C3").Value = "NSE:AARTISURF"
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "NSE:FEDERALBNK"

# --- Row 4 ---
$ws.Range("B4").Value = "NSE:LOTUSEYE"
$ws.Range("C4").Value = "NSE:AEROFLEX"
$ws.Range("F4").ClearContents()

# --- Row 5 ---
$ws.Range("B5").Value = "NSE:OIL"
$ws.Range("C5").Value = "NSE:ALKEM"

# --- Row 6 ---
$ws.Range("B6").Value = "NSE:SAKSOFT"
$ws.Range("C6").Value = "NSE:ANANTRAJ"

# --- Row 7 ---
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "NSE:CAPTRUST"

# --- Row 8 ---
$ws.Range("B8").ClearContents()
$ws.Range("C8").Value = "NSE:CEATLTD"

# --- Row 9 ---
$ws.Range("B9").ClearContents()
$ws.Range("C9").Value = "NSE:DIAMINESQ"

# --- Row 10 ---
$ws.Range("B10").ClearContents()
$ws.Range("C10").Value = "NSE:DICIND"

# --- Row 11 ---
$ws.Range("B11").ClearContents()
$ws.Range("C11").Value = "NSE:EMAMILTD"

# --- Row 12 ---
$ws.Range("B12").ClearContents()
$ws.Range("C12").Value = "NSE:GODREJIND"

# --- Row 13 ---
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "NSE:GRMOVER"

# --- Row 14 ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "NSE:HGINFRA"

# --- Row 15 ---
$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = "NSE:HINDCOPPER"

# --- Row 16 ---
$ws.Range("C16").Value = "NSE:IEL"

# --- Row 17 ---
$ws.Range("C17").Value = "NSE:INDIANCARD"

# --- Row 18 ---
$ws.Range("C18").Value = "NSE:JSWHL"

# --- Row 19 ---
$ws.Range("C19").Value = "NSE:JUBLPHARMA"

# --- Row 20 ---
$ws.Range("C20").Value = "NSE:KOHINOOR"

# --- New rows 21-27 ---
$newRows = @(
    @{ Row = 21; A = 19; C = "NSE:KOPRAN" },
    @{ Row = 22; A = 20; C = "NSE:LLOYDSENGG" },
    @{ Row = 23; A = 21; C = "NSE:PIONEEREMB" },
    @{ Row = 24; A = 22; C = "NSE:POONAWALLA" },
    @{ Row = 25; A = 23; C = "NSE:RAMRAT" },
    @{ Row = 26; A = 24; C = "NSE:ROSSELLIND" },
    @{ Row = 27; A = 25; C = "NSE:RUSHIL" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    # Copy the formatting of column A from the row above chain (row 20 has the base style)
    $ws.Range("A20").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("C$rowNum").Value = $r.C
}

$excel.CutCopyMode = $false
